# Weekly refresh of the Orégano (Vega Central Mapocho de Santiago) dataset:
# a new weekly observation is inserted as row 29 (pushing the existing
# rows 29-93 down to 30-94), and the sheet's used range grows by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 29, shifting rows 29..93 -> 30..94.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Cells.Item(29, 1).Value = 9
$ws.Cells.Item(29, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(29, 3).Value = 'Metropolitana'
$ws.Cells.Item(29, 4).Value = 44979
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100112029
$ws.Cells.Item(29, 7).Value = 'Orégano'
$ws.Cells.Item(29, 8).Value = 'Sin especificar'
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 16
$ws.Cells.Item(29, 11).Value = 18000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 18000
$ws.Cells.Item(29, 14).Value = '$/docena de atados'
$ws.Cells.Item(29, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 16).Value = 6000
$ws.Cells.Item(29, 17).Value = 3
$ws.Cells.Item(29, 18).Value = 'Hortaliza'
